$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the post entry that was in row 626 ("花より団子おおお（違う）").
# Deleting the entire row shifts rows 627:766 up to 626:765, matching the
# rest of the table, and shrinks the used range to A1:C765.
$ws.Rows.Item(626).Delete()
